$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was D-Wave Quantum Inc. / QBTS -> becomes International Business Machines / IBM
$ws.Range("B2").Value = "International Business Machines"
$ws.Range("C2").Value = "IBM"
$ws.Range("D2").Value = 307.13
$ws.Range("F2").Value = 5.76
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 63
$ws.Range("N2").Value = 85.83574689470727

# Row 3: was International Business Machines / IBM -> becomes D-Wave Quantum Inc. / QBTS
$ws.Range("B3").Value = "D-Wave Quantum Inc."
$ws.Range("C3").Value = "QBTS"
$ws.Range("D3").Value = 22.67
$ws.Range("F3").Value = 10.53
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 76
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 62.2
$ws.Range("N3").Value = 85.83574689470727

# Row 4 (Rigetti): only MACRO_SCORE refreshed
$ws.Range("N4").Value = 85.83574689470727

# Row 5 (IonQ): only MACRO_SCORE refreshed
$ws.Range("N5").Value = 85.83574689470727
